# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# D-column "Price" values are stored as text (they can contain two dots,
# e.g. "30.675.35"), so they are written with a leading apostrophe to stop
# Excel's automatic number coercion from mangling them (dropping trailing
# zeros, switching to scientific notation, etc). E-column values already
# contain surrounding spaces/percent signs so they stay text on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.675.35"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "'2.121.01"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +1.10%  "
$ws.Range("D5").Value = "'337.73"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").Value = "'1.011"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("D7").Value = "'0.5254"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").Value = "'0.4558"
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("D9").Value = "'54.86"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("D10").Value = "'0.09124"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").Value = "'1.175"
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").Value = "'24.50"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "'2.115.87"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").Value = "'8.160"
$ws.Range("E15").Value = "  +5.83%  "
$ws.Range("D16").Value = "'0.00001176"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("D17").Value = "'97.30"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "'0.06698"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").Value = "'19.49"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "'6.324"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "'30.752.71"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").Value = "'12.90"
$ws.Range("E24").Value = "  +4.93%  "
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").Value = "'2.374.00"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").Value = "'22.42"
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'165.00"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.566"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "'134.81"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").Value = "'1.667"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "'6.385"
$ws.Range("E34").Value = "  +3.54%  "
$ws.Range("D35").Value = "'3.945"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("D36").Value = "'10.65"
$ws.Range("E36").Value = "  +3.60%  "
$ws.Range("D37").Value = "'5.891"
$ws.Range("E37").Value = "  +7.68%  "
$ws.Range("D38").Value = "'0.02641"
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").Value = "'0.2327"
$ws.Range("E40").Value = "  +2.71%  "
$ws.Range("D41").Value = "'12.69"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").Value = "'0.6925"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "'15.20"
$ws.Range("E44").Value = "  +9.01%  "
$ws.Range("D45").Value = "'0.6498"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").Value = "'2.324"
$ws.Range("E46").Value = "  +2.98%  "
$ws.Range("D47").Value = "'0.00000000373"
$ws.Range("E47").Value = "  +22.50%  "
$ws.Range("D48").Value = "'3.697"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").Value = "'1.257"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").Value = "'83.50"
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("D51").Value = "'0.07315"
$ws.Range("E51").Value = "  +3.92%  "
